$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("s6",  "s6_IMG_3178.jpeg",  "meltpatch", "1705", "1748", "6", "3", "31",  "2"),
    @("s7",  "s7_IMG_3180.jpeg",  "meltpatch", "618",  "529",  "6", "3", "150", "2"),
    @("s8",  "s8_IMG_3176.jpeg",  "meltpatch", "460",  "1320", "6", "3", "111", "2"),
    @("s9",  "s9_IMG_3177.jpeg",  "meltpatch", "3353", "1839", "6", "3", "106", "2"),
    @("s10", "s10_IMG_3174.jpeg", "meltpatch", "3045", "2170", "6", "3", "3",   "2")
)

$startRow = 7
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $cell = $ws.Cells.Item($row, $c + 1)
        $cell.NumberFormat = "@"
        $cell.Value = $vals[$c]
        $cell.Style = "Normal"
    }
}
